# Hortaliza, Vega Modelo de Temuco - Acelga
# Insert two new weekly price rows (502 and 503) into the daily logic
# table, shifting the existing rows 502:584 down to 504:586.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing rows down by two to make room for the new entries.
$ws.Rows("502:503").Insert()

# --- New row 502 ---------------------------------------------------------
$ws.Cells.Item(502, 1).Value = 10
$ws.Cells.Item(502, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(502, 3).Value = "La Araucanía"
$ws.Cells.Item(502, 4).Value = 45218
$ws.Cells.Item(502, 5).Value = 9
$ws.Cells.Item(502, 6).Value = 100112009
$ws.Cells.Item(502, 7).Value = "Acelga"
$ws.Cells.Item(502, 8).Value = "Sin especificar"
$ws.Cells.Item(502, 9).Value = "Primera"
$ws.Cells.Item(502, 10).Value = 65
$ws.Cells.Item(502, 11).Value = 8000
$ws.Cells.Item(502, 12).Value = 8000
$ws.Cells.Item(502, 13).Value = 8000
$ws.Cells.Item(502, 14).Value = "`$/docena de atados (12 kilos)"
$ws.Cells.Item(502, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(502, 16).Value = 667
$ws.Cells.Item(502, 17).Value = 12
$ws.Cells.Item(502, 18).Value = "Hortaliza"

# --- New row 503 ---------------------------------------------------------
$ws.Cells.Item(503, 1).Value = 10
$ws.Cells.Item(503, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(503, 3).Value = "La Araucanía"
$ws.Cells.Item(503, 4).Value = 45218
$ws.Cells.Item(503, 5).Value = 9
$ws.Cells.Item(503, 6).Value = 100112009
$ws.Cells.Item(503, 7).Value = "Acelga"
$ws.Cells.Item(503, 8).Value = "Sin especificar"
$ws.Cells.Item(503, 9).Value = "Primera"
$ws.Cells.Item(503, 10).Value = 65
$ws.Cells.Item(503, 11).Value = 7000
$ws.Cells.Item(503, 12).Value = 7000
$ws.Cells.Item(503, 13).Value = 7000
$ws.Cells.Item(503, 14).Value = "`$/docena de atados (12 kilos)"
$ws.Cells.Item(503, 15).Value = "Región del Maule"
$ws.Cells.Item(503, 16).Value = 583
$ws.Cells.Item(503, 17).Value = 12
$ws.Cells.Item(503, 18).Value = "Hortaliza"
